# Append new log rows (445-453) to Sheet1, matching the source diff.
# Columns: A=staff_name, B=id, C=start_time, D=leaving_time, E=reason,
#          F=department, G=date, H=btn_id
# All cells in this sheet are stored as literal text (no numeric/date
# auto-typing), so columns B (numeric-looking ids) and G (date-looking
# strings) must be pre-formatted as Text before the value is entered -
# otherwise Excel auto-converts them to a number / date serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("NGUYỄN TUẤN ANH",    "224057", "16:30", "16:30", "TEST REQUEST", "RD", "2024-05-20", "btn_14"),
    @("PHẠM THỊ PHƯƠNG",    "172684", "16:30", "19:30", "TEST REQUEST", "RD", "2024-05-20", "btn_5"),
    @("TRƯƠNG VĂN MINH",    "203591", "16:30", "19:30", "TEST REQUEST", "RD", "2024-05-20", "btn_9"),
    @("LÊ MINH THẮNG",      "223906", "16:30", "19:30", "TEST REQUEST", "RD", "2024-05-20", "btn_12"),
    @("NGUYỄN HOÀNG VIỆT",  "172759", "16:30", "19:30", "TEST REQUEST", "RD", "2024-05-20", "btn_6"),
    @("TRẦN VĂN LƯU",       "234102", "16:30", "19:30", "TEST REQUEST", "RD", "2024-05-20", "btn_15"),
    @("LÊ HUỲNH ANH KHOA",  "234168", "16:30", "19:30", "B/T",          "RD", "2024-05-20", "btn_16"),
    @("NGUYỄN QUANG QUÍ",   "203638", "16:30", "19:30", "TEST REQUEST", "RD", "2024-05-20", "btn_10"),
    @("NGUYỄN TUẤN ANH",    "224057", "16:30", "14:30", "TEST REQUEST", "RD", "2024-05-21", "btn_14")
)

$startRow = 445
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]

    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $data[1]

    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]

    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = $data[6]

    $ws.Cells.Item($r, 8).Value = $data[7]
}
